$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Translate header row from Korean to English
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "PhoneNumber"
$ws.Range("C1").Value = "ID"
$ws.Range("D1").Value = "Status"

# Update selection to reflect the new active cell in the frozen-pane view
$ws.Range("A2").Select()
